$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting aHR..Std Error (and the
# matching data) one column to the right.
$ws.Columns("B:B").Insert()

# Header for the new column, styled like the other headers.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Predictor"

# Fill the new column with the same predictor labels as column A, as
# plain (unstyled) text cells.
$ws.Range("B2:B6").ClearFormats()
$ws.Range("B2").Value = $ws.Range("A2").Value()
$ws.Range("B3").Value = $ws.Range("A3").Value()
$ws.Range("B4").Value = $ws.Range("A4").Value()
$ws.Range("B5").Value = $ws.Range("A5").Value()
$ws.Range("B6").Value = $ws.Range("A6").Value()

Write-Output "done"
